# Refresh the cryptos price table (GitHub Actions scheduled data pull).
# Updates the "Price" (D) and "Volume(1h)" (E) columns for the rows whose
# figures moved, and swaps the Cosmos/Monero rows (25 <-> 26) which
# reordered in the source ranking.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column holds plain-text figures (e.g. "233.97", "0.4666")
# that are NOT real numbers in the sheet (no numeric formatting is applied
# to them). Several of the new values look like ordinary decimals, so
# Excel's type-inference on Range.Value would silently coerce them into
# numbers (dropping formatting like "12.60" -> 12.6). Force those specific
# cells to Text format first so the literal strings are preserved exactly.
$textRefs = @("D5", "D7", "D8", "D9", "D10", "D11", "D12", "D14", "D15", "D16", "D19", "D20", "D22", "D24", "D25", "D26", "D27", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D46", "D48", "D49", "D50", "D51")
foreach ($r in $textRefs) { $ws.Range($r).NumberFormat = "@" }

$ws.Range("D2").Value = "30.142.60"
$ws.Range("E2").Value = "  -1.85%  "
$ws.Range("D3").Value = "1.857.30"
$ws.Range("E3").Value = "  -3.42%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").Value = "233.97"
$ws.Range("E5").Value = "  -3.09%  "
$ws.Range("E6").Value = "  +0.15%  "
$ws.Range("D7").Value = "0.4666"
$ws.Range("E7").Value = "  -2.47%  "
$ws.Range("D8").Value = "0.2815"
$ws.Range("E8").Value = "  -2.66%  "
$ws.Range("D9").Value = "0.06566"
$ws.Range("E9").Value = "  -3.04%  "
$ws.Range("D10").Value = "19.95"
$ws.Range("E10").Value = "  +1.64%  "
$ws.Range("D11").Value = "0.07831"
$ws.Range("E11").Value = "  +0.56%  "
$ws.Range("D12").Value = "96.69"
$ws.Range("D13").Value = "1.865.21"
$ws.Range("E13").Value = "  -3.07%  "
$ws.Range("D14").Value = "5.109"
$ws.Range("E14").Value = "  -3.16%  "
$ws.Range("D15").Value = "0.6652"
$ws.Range("E15").Value = "  -2.18%  "
$ws.Range("D16").Value = "282.69"
$ws.Range("E16").Value = "  -3.69%  "
$ws.Range("D17").Value = "30.183.38"
$ws.Range("E17").Value = "  -1.76%  "
$ws.Range("E18").Value = "  +0.12%  "
$ws.Range("D19").Value = "5.430"
$ws.Range("E19").Value = "  -1.21%  "
$ws.Range("D20").Value = "12.60"
$ws.Range("D21").Value = "2.110.97"
$ws.Range("E21").Value = "  -3.27%  "
$ws.Range("D22").Value = "0.000007237"
$ws.Range("E22").Value = "  -4.45%  "
$ws.Range("E23").Value = "  +0.11%  "
$ws.Range("D24").Value = "6.138"
$ws.Range("E24").Value = "  -3.85%  "
$ws.Range("B25").Value = "Monero"
$ws.Range("C25").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D25").Value = "167.70"
$ws.Range("E25").Value = "  -0.42%  "
$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D26").Value = "9.325"
$ws.Range("E26").Value = "  -1.98%  "
$ws.Range("D27").Value = "18.89"
$ws.Range("E27").Value = "  -4.78%  "
$ws.Range("E28").Value = "  -9.37%  "
$ws.Range("D29").Value = "1.335"
$ws.Range("E29").Value = "  -4.19%  "
$ws.Range("D30").Value = "0.09569"
$ws.Range("E30").Value = "  -4.92%  "
$ws.Range("D31").Value = "4.406"
$ws.Range("E31").Value = "  -4.66%  "
$ws.Range("D32").Value = "1.471"
$ws.Range("E32").Value = "  -3.67%  "
$ws.Range("D33").Value = "4.102"
$ws.Range("E33").Value = "  -4.84%  "
$ws.Range("D34").Value = "0.04667"
$ws.Range("E34").Value = "  -2.70%  "
$ws.Range("D35").Value = "0.7005"
$ws.Range("E35").Value = "  -4.68%  "
$ws.Range("E36").Value = "  -2.24%  "
$ws.Range("D37").Value = "1.0000"
$ws.Range("E37").Value = "  +0.23%  "
$ws.Range("D38").Value = "2.700"
$ws.Range("E38").Value = "  -0.51%  "
$ws.Range("D39").Value = "0.01851"
$ws.Range("E39").Value = "  -4.51%  "
$ws.Range("D40").Value = "6.333"
$ws.Range("E40").Value = "  -1.27%  "
$ws.Range("D41").Value = "2.509"
$ws.Range("E41").Value = "  -4.40%  "
$ws.Range("D42").Value = "71.97"
$ws.Range("E42").Value = "  -4.13%  "
$ws.Range("D43").Value = "0.8533"
$ws.Range("E43").Value = "  -1.41%  "
$ws.Range("D44").Value = "1.926"
$ws.Range("E44").Value = "  -3.60%  "
$ws.Range("E45").Value = "  +0.11%  "
$ws.Range("D46").Value = "0.4159"
$ws.Range("E46").Value = "  -3.72%  "
$ws.Range("E47").Value = "  -2.39%  "
$ws.Range("D48").Value = "993.99"
$ws.Range("E48").Value = "  +1.19%  "
$ws.Range("D49").Value = "7.230"
$ws.Range("E49").Value = "  -4.05%  "
$ws.Range("D50").Value = "9.175"
$ws.Range("E50").Value = "  +2.40%  "
$ws.Range("D51").Value = "33.92"
$ws.Range("E51").Value = "  -2.93%  "
